# "contingencies with rene fine"
# Update existing rows 8-15 (line7/line8 take over what used to be the first two
# "extr" rows, and every row below shifts its label by two: extr1..extr8 move
# down to rows 10-17) and append two brand-new rows (16-17) for extr7/extr8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: now "line7" ---
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- Row 9: now "line8" ---
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Row 10: now "extr1" ---
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# --- Row 11: now "extr2" ---
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# --- Row 12: now "extr3" ---
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

# --- Row 13: now "extr4" ---
$ws.Range("B13").Value = "extr4"
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

# --- Row 14: now "extr5" ---
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# --- Row 15: now "extr6" ---
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- New row 16: "extr7" (copy formatting from row 15 first) ---
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# --- New row 17: "extr8" ---
$ws.Range("A15:E15").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

$excel.CutCopyMode = $false
